$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$nl = [char]10

# ---- Cell values (order matters: it determines sharedStrings index order) ----
$ws.Cells.Item(3,1).Value = "Apache common io"
$ws.Cells.Item(3,2).Value = "read a file"
$ws.Cells.Item(3,3).Value = "import org.apache.commona.io.FileUtils;" + $nl + "String fileContent=FileUtils.readFileToString(file);"

$ws.Cells.Item(4,1).Value = "java"
$ws.Cells.Item(4,2).Value = "How to run a executable jar"

$ws.Cells.Item(5,1).Value = "Maven"
$ws.Cells.Item(5,2).Value = "build a java project"
$ws.Cells.Item(5,3).Value = "`$ mvn archetype:generate " + $nl + "        -DgroupId=com.mycompany.app" + $nl + "        -DartifactId=my-app " + $nl + "        -DarchetypeArtifactId=maven-archetype-quickstart " + $nl + "        -DinteractiveMode=false"

$ws.Cells.Item(6,1).Value = "Maven"
$ws.Cells.Item(6,2).Value = "import project to eclipse"
$ws.Cells.Item(6,3).Value = "1. Generate necessary config file for eclipse:" + $nl + "`$mvn eclipse:eclipse -Dwtpversion=2.0" + $nl + "2. Imports it into Eclipse IDE" + $nl + "File -> Import… -> General -> Existing Projects into workspace"

$ws.Cells.Item(7,1).Value = "Maven"
$ws.Cells.Item(7,2).Value = "build a web project"
$ws.Cells.Item(7,3).Value = "`$ mvn archetype:generate " + $nl + "        -DgroupId={project-packaging} " + $nl + "        -DartifactId={project-name} " + $nl + "        -DarchetypeArtifactId=maven-archetype-webapp " + $nl + "        -DinteractiveMode=false"

$ws.Cells.Item(8,1).Value = "Tomcat"
$ws.Cells.Item(8,2).Value = "basic command"
$ws.Cells.Item(8,3).Value = "`$/etc/init.d/`$ sudo ./tomcat6 start" + $nl + "`$/etc/init.d/`$ sudo ./tomcat6 stop" + $nl + "`$/etc/init.d/`$ sudo ./tomcat6 restart "

# This one is added last on purpose, to match the target sharedStrings order (index 21)
$ws.Cells.Item(4,3).Value = "1. java -jar *.jar {argus} : 直接針對該JAR調用MAIN" + $nl + "2. java -classpath *.jar com.*.*.* : 曲折地把該JAR加為CP, 然後用CP概念調用任一包含STATIC方法的類"

Write-Host "Done values"
